$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the "site_id" row (row 2), pushing every
# other field row down by 3 (row 3 -> row 6, ..., row 13 -> row 16).
$ws.Rows("3:5").Insert()

# Copy the formatting of the "site_id" row (row 2) down into the three
# freshly-inserted rows so they pick up the same style indices (bold/
# bordered field-name cell, bordered description/sample cells, etc.)
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E5").PasteSpecial(-4122)

# Row 3: site_name. The "Sample" column is intentionally left blank for
# this field.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "site_name"
$ws.Range("C3").Value = "varchar"
$ws.Range("D3").Value = "name of the measurement site"

# Row 4 + 5 Field column: state_id / state_name
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "state_id"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "state_name"

# Row 4 + 5 Description column: state code / state name
$ws.Range("D4").Value = "state code"
$ws.Range("D5").Value = "state name"

# Row 4 Type column: smallint, row 5 Type column: varchar
$ws.Range("C4").Value = "smallint"
$ws.Range("C5").Value = "varchar"

# Row 4 + 5 Sample column: 6 / California
$ws.Range("E4").Value = "6"
$ws.Range("E5").Value = "California"

# Renumber the "#" column for the rows that shifted down (previously 3..12,
# now located at rows 6..16, values should read 5..15).
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15

# Drop the stale selection left over from the original sheet view so it
# doesn't keep pointing at the old A1:E13 extent.
$ws.Range("A1").Select()
